$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the entire used range first, then write the new (smaller) table.
$ws.Cells.Clear()

# Headers
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Data row
$ws.Range("A2").Value = 81
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 8267
$ws.Range("D2").Value = 0.1435472965240479
